$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.126994
$ws.Range("H2").Value = 0.380982
$ws.Range("I2").Value = 0.1347536544371063
$ws.Range("J2").Value = 0.1347536544371063
$ws.Range("M2").Value = 1.089622333333333
$ws.Range("N2").Value = 3.268867
$ws.Range("O2").Value = 0.09062709179941439
$ws.Range("P2").Value = 0.09062709179941439
$ws.Range("Q2").Value = 0.1383754985993333
$ws.Range("R2").Value = 1.245379487394
$ws.Range("S2").Value = 0.0122123318109782
$ws.Range("T2").Value = 0.0122123318109782
$ws.Range("G3").Value = 0.126994
$ws.Range("H3").Value = 0.380982
$ws.Range("I3").Value = 0.1347536544371063
$ws.Range("J3").Value = 0.1347536544371063
$ws.Range("O3").Value = 0.1560803133424032
$ws.Range("P3").Value = 0.1560803133424032
$ws.Range("Q3").Value = 0.2383138502126667
$ws.Range("R3").Value = 2.144824651914
$ws.Range("S3").Value = 0.02103239260857748
$ws.Range("T3").Value = 0.02103239260857748
$ws.Range("G4").Value = 0.126994
$ws.Range("H4").Value = 0.380982
$ws.Range("I4").Value = 0.1347536544371063
$ws.Range("J4").Value = 0.1347536544371063
$ws.Range("M4").Value = 0.6928603333333333
$ws.Range("N4").Value = 2.078581
$ws.Range("O4").Value = 0.05762723019918477
$ws.Range("P4").Value = 0.05762723019918477
$ws.Range("Q4").Value = 0.08798910517133332
$ws.Range("R4").Value = 0.7919019465419999
$ws.Range("S4").Value = 0.007765479864428522
$ws.Range("T4").Value = 0.007765479864428522
$ws.Range("G5").Value = 0.126994
$ws.Range("H5").Value = 0.380982
$ws.Range("I5").Value = 0.1347536544371063
$ws.Range("J5").Value = 0.1347536544371063
$ws.Range("M5").Value = 8.364082999999999
$ws.Range("N5").Value = 25.092249
$ws.Range("O5").Value = 0.6956653646589975
$ws.Range("P5").Value = 0.6956653646589976
$ws.Range("Q5").Value = 1.062188356502
$ws.Range("R5").Value = 9.559695208517999
$ws.Range("S5").Value = 0.09374345015312212
$ws.Range("T5").Value = 0.09374345015312213
$ws.Range("I6").Value = 0.1440694272310034
$ws.Range("J6").Value = 0.1440694272310034
$ws.Range("M6").Value = 1.089622333333333
$ws.Range("N6").Value = 3.268867
$ws.Range("O6").Value = 0.09062709179941439
$ws.Range("P6").Value = 0.09062709179941439
$ws.Range("Q6").Value = 0.1479416562711111
$ws.Range("R6").Value = 1.33147490644
$ws.Range("S6").Value = 0.0130565932071532
$ws.Range("T6").Value = 0.0130565932071532
$ws.Range("I7").Value = 0.1440694272310034
$ws.Range("J7").Value = 0.1440694272310034
$ws.Range("O7").Value = 0.1560803133424032
$ws.Range("P7").Value = 0.1560803133424032
$ws.Range("S7").Value = 0.02248640134527558
$ws.Range("T7").Value = 0.02248640134527558
$ws.Range("I8").Value = 0.1440694272310034
$ws.Range("J8").Value = 0.1440694272310034
$ws.Range("M8").Value = 0.6928603333333333
$ws.Range("N8").Value = 2.078581
$ws.Range("O8").Value = 0.05762723019918477
$ws.Range("P8").Value = 0.05762723019918477
$ws.Range("Q8").Value = 0.09407195699111111
$ws.Range("R8").Value = 0.8466476129199999
$ws.Range("S8").Value = 0.008302322047705734
$ws.Range("T8").Value = 0.008302322047705734
$ws.Range("I9").Value = 0.1440694272310034
$ws.Range("J9").Value = 0.1440694272310034
$ws.Range("M9").Value = 8.364082999999999
$ws.Range("N9").Value = 25.092249
$ws.Range("O9").Value = 0.6956653646589975
$ws.Range("P9").Value = 0.6956653646589976
$ws.Range("Q9").Value = 1.135619429186667
$ws.Range("R9").Value = 10.22057486268
$ws.Range("S9").Value = 0.1002241106308689
$ws.Range("T9").Value = 0.1002241106308689
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2408096666666667
$ws.Range("H10").Value = 0.722429
$ws.Range("I10").Value = 0.2555237460597726
$ws.Range("J10").Value = 0.2555237460597726
$ws.Range("M10").Value = 1.089622333333333
$ws.Range("N10").Value = 3.268867
$ws.Range("O10").Value = 0.09062709179941439
$ws.Range("P10").Value = 0.09062709179941439
$ws.Range("Q10").Value = 0.2623915908825556
$ws.Range("R10").Value = 2.361524317943
$ws.Range("S10").Value = 0.02315737399108926
$ws.Range("T10").Value = 0.02315737399108926
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2408096666666667
$ws.Range("H11").Value = 0.722429
$ws.Range("I11").Value = 0.2555237460597726
$ws.Range("J11").Value = 0.2555237460597726
$ws.Range("O11").Value = 0.1560803133424032
$ws.Range("P11").Value = 0.1560803133424032
$ws.Range("Q11").Value = 0.4518975607647778
$ws.Range("R11").Value = 4.067078046883
$ws.Range("S11").Value = 0.03988222635143399
$ws.Range("T11").Value = 0.03988222635143399
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2408096666666667
$ws.Range("H12").Value = 0.722429
$ws.Range("I12").Value = 0.2555237460597726
$ws.Range("J12").Value = 0.2555237460597726
$ws.Range("M12").Value = 0.6928603333333333
$ws.Range("N12").Value = 2.078581
$ws.Range("O12").Value = 0.05762723019918477
$ws.Range("P12").Value = 0.05762723019918477
$ws.Range("Q12").Value = 0.1668474659165556
$ws.Range("R12").Value = 1.501627193249
$ws.Range("S12").Value = 0.01472512573554455
$ws.Range("T12").Value = 0.01472512573554455
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2408096666666667
$ws.Range("H13").Value = 0.722429
$ws.Range("I13").Value = 0.2555237460597726
$ws.Range("J13").Value = 0.2555237460597726
$ws.Range("M13").Value = 8.364082999999999
$ws.Range("N13").Value = 25.092249
$ws.Range("O13").Value = 0.6956653646589975
$ws.Range("P13").Value = 0.6956653646589976
$ws.Range("Q13").Value = 2.014152039202333
$ws.Range("R13").Value = 18.127368352821
$ws.Range("S13").Value = 0.1777590199817048
$ws.Range("T13").Value = 0.1777590199817048
$ws.Range("G14").Value = 0.438839
$ws.Range("H14").Value = 1.316517
$ws.Range("I14").Value = 0.4656531722721176
$ws.Range("J14").Value = 0.4656531722721176
$ws.Range("M14").Value = 1.089622333333333
$ws.Range("N14").Value = 3.268867
$ws.Range("O14").Value = 0.09062709179941439
$ws.Range("P14").Value = 0.09062709179941439
$ws.Range("Q14").Value = 0.4781687751376667
$ws.Range("R14").Value = 4.303518976239
$ws.Range("S14").Value = 0.04220079279019372
$ws.Range("T14").Value = 0.04220079279019372
$ws.Range("G15").Value = 0.438839
$ws.Range("H15").Value = 1.316517
$ws.Range("I15").Value = 0.4656531722721176
$ws.Range("J15").Value = 0.4656531722721176
$ws.Range("O15").Value = 0.1560803133424032
$ws.Range("P15").Value = 0.1560803133424032
$ws.Range("Q15").Value = 0.8235145889843333
$ws.Range("R15").Value = 7.411631300859
$ws.Range("S15").Value = 0.07267929303711619
$ws.Range("T15").Value = 0.07267929303711619
$ws.Range("G16").Value = 0.438839
$ws.Range("H16").Value = 1.316517
$ws.Range("I16").Value = 0.4656531722721176
$ws.Range("J16").Value = 0.4656531722721176
$ws.Range("M16").Value = 0.6928603333333333
$ws.Range("N16").Value = 2.078581
$ws.Range("O16").Value = 0.05762723019918477
$ws.Range("P16").Value = 0.05762723019918477
$ws.Range("Q16").Value = 0.3040541358196666
$ws.Range("R16").Value = 2.736487222376999
$ws.Range("S16").Value = 0.02683430255150596
$ws.Range("T16").Value = 0.02683430255150596
$ws.Range("G17").Value = 0.438839
$ws.Range("H17").Value = 1.316517
$ws.Range("I17").Value = 0.4656531722721176
$ws.Range("J17").Value = 0.4656531722721176
$ws.Range("M17").Value = 8.364082999999999
$ws.Range("N17").Value = 25.092249
$ws.Range("O17").Value = 0.6956653646589975
$ws.Range("P17").Value = 0.6956653646589976
$ws.Range("Q17").Value = 3.670485819636999
$ws.Range("R17").Value = 33.034372376733
$ws.Range("S17").Value = 0.3239387838933017
$ws.Range("T17").Value = 0.3239387838933018
